# Update weekly Fruta/Hortaliza price records for
# "Femacal de La Calera - Caqui" sheet.
# The underlying dataset rows were re-shuffled across dates while keeping
# the same overall multiset of values; below are the per-row corrections
# needed to go from the old layout to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44307
$ws.Range("M2").Value = 40

# Row 3
$ws.Range("D3").Value = 44301
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("S3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44314
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 47
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 9000
$ws.Range("S4").Value = 900

# Row 5
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 58
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 9000
$ws.Range("S6").Value = 900

# Row 7
$ws.Range("D7").Value = 44333
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("S7").Value = 800

# Row 11
$ws.Range("D11").Value = 44319
$ws.Range("M11").Value = 68
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("R11").Value = "Provincia de Quillota"
$ws.Range("S11").Value = 1000

# Row 12
$ws.Range("D12").Value = 44319
$ws.Range("M12").Value = 57
$ws.Range("R12").Value = "Provincia de Quillota"

# Row 13
$ws.Range("D13").Value = 44326
$ws.Range("M13").Value = 65

# Row 14
$ws.Range("D14").Value = 44326
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 67
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("S14").Value = 800

# Row 15
$ws.Range("D15").Value = 44315
$ws.Range("M15").Value = 45

# Row 17
$ws.Range("D17").Value = 44312
$ws.Range("M17").Value = 48

# Row 18
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("S18").Value = 1000

# Row 19
$ws.Range("D19").Value = 44323
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 9000
$ws.Range("O19").Value = 9000
$ws.Range("P19").Value = 9000
$ws.Range("S19").Value = 900

# Row 20
$ws.Range("D20").Value = 44306
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 45
$ws.Range("N20").Value = 10000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 10000
$ws.Range("S20").Value = 1000

# Row 21
$ws.Range("D21").Value = 44328
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 800

# Row 22
$ws.Range("D22").Value = 44328
$ws.Range("N22").Value = 7000
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 7000
$ws.Range("S22").Value = 700

# Row 23
$ws.Range("D23").Value = 44322
$ws.Range("M23").Value = 56

# Row 24
$ws.Range("D24").Value = 44322
$ws.Range("M24").Value = 40

# Row 25
$ws.Range("D25").Value = 44321
$ws.Range("M25").Value = 58
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("S25").Value = 900

# Row 26
$ws.Range("D26").Value = 44308

# Row 27
$ws.Range("D27").Value = 44308
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 48
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 8000
$ws.Range("S27").Value = 800

# Row 28
$ws.Range("D28").Value = 44329
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 56
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 9000
$ws.Range("P28").Value = 9000
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 900

# Row 29
$ws.Range("D29").Value = 44329
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 8000
$ws.Range("O29").Value = 8000
$ws.Range("P29").Value = 8000
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 800

# Row 30
$ws.Range("D30").Value = 44302
$ws.Range("M30").Value = 45
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("S30").Value = 1000
